$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 86
$ws.Range("I2").Value = 293
$ws.Range("J2").Value = 1099
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 276
$ws.Range("M2").Value = 14
$ws.Range("N2").Value = 184
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 14
$ws.Range("S2").Value = 147
$ws.Range("T2").Value = 182
$ws.Range("U2").Value = 17
$ws.Range("V2").Value = 1719
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1649
$ws.Range("Y2").Value = 3
$ws.Range("AA2").Value = 9
